$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 2.15
$ws.Range("H2").Value = 3.2
$ws.Range("K2").Value = 2
$ws.Range("L2").Value = 4.33
$ws.Range("O2").Value = 1.36
$ws.Range("P2").Value = 3
$ws.Range("Q2").Value = 2.2
$ws.Range("R2").Value = 1.65
$ws.Range("S2").Value = 1.5
$ws.Range("T2").Value = 2.5
$ws.Range("U2").Value = 1.91
$ws.Range("V2").Value = 1.8
$ws.Range("W2").Value = 6.5
$ws.Range("Y2").Value = 9.5
$ws.Range("AB2").Value = 34
$ws.Range("AC2").Value = 8
$ws.Range("AE2").Value = 17
$ws.Range("AG2").Value = 9
$ws.Range("AO2").Value = 26
$ws.Range("AR2").Value = 2.5
$ws.Range("AT2").Value = 67
$ws.Range("AW2").Value = 34
$ws.Range("BA2").Value = 351
